$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.681896924972534
$ws.Range("B1").Value = 7.013472557067871
$ws.Range("C1").Value = 7.132493495941162
$ws.Range("D1").Value = 2.379313707351685
$ws.Range("E1").Value = 1.46185314655304
